$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at row 51, pushing existing rows (old 51-98) down to (53-100)
$ws.Rows("51:52").Insert()

# Copy the style (incl. number format) used on column D date cells to the new rows
$ws.Range("D53").Copy()
$ws.Range("D51:D52").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 51 data
$ws.Cells.Item(51, 1).Value = 11
$ws.Cells.Item(51, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(51, 3).Value = "Bíobío"
$ws.Cells.Item(51, 4).Value = 44601
$ws.Cells.Item(51, 5).Value = 8
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100102
$ws.Cells.Item(51, 8).Value = "Cítricos"
$ws.Cells.Item(51, 9).Value = 100102004
$ws.Cells.Item(51, 10).Value = "Mandarina"
$ws.Cells.Item(51, 11).Value = "Murcott"
$ws.Cells.Item(51, 12).Value = "Primera"
$ws.Cells.Item(51, 13).Value = 100
$ws.Cells.Item(51, 14).Value = 9000
$ws.Cells.Item(51, 15).Value = 10000
$ws.Cells.Item(51, 16).Value = 9500
$ws.Cells.Item(51, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(51, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(51, 19).Value = 528
$ws.Cells.Item(51, 20).Value = 18

# Row 52 data
$ws.Cells.Item(52, 1).Value = 11
$ws.Cells.Item(52, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(52, 3).Value = "Bíobío"
$ws.Cells.Item(52, 4).Value = 44601
$ws.Cells.Item(52, 5).Value = 8
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100102
$ws.Cells.Item(52, 8).Value = "Cítricos"
$ws.Cells.Item(52, 9).Value = 100102004
$ws.Cells.Item(52, 10).Value = "Mandarina"
$ws.Cells.Item(52, 11).Value = "Murcott"
$ws.Cells.Item(52, 12).Value = "Segunda"
$ws.Cells.Item(52, 13).Value = 50
$ws.Cells.Item(52, 14).Value = 8000
$ws.Cells.Item(52, 15).Value = 8000
$ws.Cells.Item(52, 16).Value = 8000
$ws.Cells.Item(52, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(52, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(52, 19).Value = 444
$ws.Cells.Item(52, 20).Value = 18
